$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal de travail")

# --- Row 53: Temps [h] 6 -> 7 ---
$ws.Cells.Item(53, 3).Value = 7

# --- Row 54: new journal entry (date, type, temps, travail effectue) ---
$ws.Cells.Item(54, 1).Value = (Get-Date -Year 2023 -Month 6 -Day 19).Date
$ws.Cells.Item(54, 2).Value = "Implémentation"
$ws.Cells.Item(54, 3).Value = 3
$ws.Cells.Item(54, 4).Value = "Correctifs backend + gestion de l'édition du dashboard sur frontend"

# --- Scroll the sheet view so row 36 is the top-left visible row ---
$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 36

# Recalculate so the totals row (row 62, column C SUM formula) reflects new data
$excel.Calculate()
